$d = $word.ActiveDocument

$replacements = @(
    @("36×72=2592", "85×80=6800"),
    @("94×83=7802", "12×98=1176"),
    @("12×15=180", "33×75=2475"),
    @("71×35=2485", "87×47=4089"),
    @("47×77=3619", "88×14=1232"),
    @("46×30=1380", "51×93=4743"),
    @("13×96=1248", "49×70=3430"),
    @("88×23=2024", "18×16=288"),
    @("79×51=4029", "69×78=5382"),
    @("91×97=8827", "43×70=3010"),
    @("15×31=465", "81×60=4860"),
    @("26×75=1950", "73×13=949"),
    @("32×11=352", "40×72=2880"),
    @("46×40=1840", "68×11=748"),
    @("21×29=609", "88×20=1760"),
    @("31×63=1953", "72×44=3168"),
    @("82×80=6560", "45×28=1260"),
    @("29×76=2204", "47×49=2303"),
    @("80×87=6960", "44×64=2816"),
    @("29×93=2697", "12×92=1104"),
    @("64×48=3072", "32×98=3136"),
    @("19×83=1577", "87×25=2175"),
    @("69×45=3105", "72×22=1584"),
    @("22×55=1210", "63×57=3591"),
    @("85×38=3230", "17×56=952")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
